$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, pushing existing row 37 (and below) down to row 38.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Range("A37").Value() = 10
$ws.Range("B37").Value() = "Vega Modelo de Temuco"
$ws.Range("C37").Value() = "La Araucanía"
$ws.Range("D37").Value() = 44959
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("E37").Value() = 9
$ws.Range("F37").Value() = 100112022
$ws.Range("G37").Value() = "Arveja Verde"
$ws.Range("H37").Value() = "Sin especificar"
$ws.Range("I37").Value() = "Primera"
$ws.Range("J37").Value() = 50
$ws.Range("K37").Value() = 35000
$ws.Range("L37").Value() = 35000
$ws.Range("M37").Value() = 35000
$ws.Range("N37").Value() = "$/saco 25 kilos"
$ws.Range("O37").Value() = "Región de La Araucanía"
$ws.Range("P37").Value() = 1400
$ws.Range("Q37").Value() = 25
$ws.Range("R37").Value() = "Hortaliza"
